$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 1356.4
$ws.Range("I2").Value2 = 1320.5
$ws.Range("K2").Value2 = 1320.5
$ws.Range("M2").Value2 = -1207.5
$ws.Range("H80").Value2 = 461.57144
$ws.Range("I80").Value2 = 722.8333
$ws.Range("K80").Value2 = 2168.4999
$ws.Range("M80").Value2 = -1170.4999
$ws.Range("H83").Value2 = 461.57144
$ws.Range("I83").Value2 = 722.8333
$ws.Range("K83").Value2 = 6505.4997
$ws.Range("M83").Value2 = -1513.4997
$ws.Range("H88").Value2 = 1166.6666
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 1166.6666
$ws.Range("K88").Value2 = 0
$ws.Range("L88").Value2 = 1166.6666
$ws.Range("N88").Value2 = -1978.6666
$ws.Range("H91").Value2 = 1166.6666
$ws.Range("I91").Value2 = 0
$ws.Range("J91").Value2 = 1166.6666
$ws.Range("K91").Value2 = 0
$ws.Range("L91").Value2 = 1166.6666
$ws.Range("N91").Value2 = -3974.6666
$ws.Range("H97").Value2 = 1027
$ws.Range("J97").Value2 = 1027
$ws.Range("L97").Value2 = 3081
$ws.Range("N97").Value2 = -4073
$ws.Range("H116").Value2 = 3655.625
$ws.Range("J116").Value2 = 4124.3335
$ws.Range("L116").Value2 = 4124.3335
$ws.Range("N116").Value2 = -11008.3335
$ws.Range("H135").Value2 = 850.2692
$ws.Range("I135").Value2 = 754.5
$ws.Range("K135").Value2 = 6790.5
$ws.Range("M135").Value2 = -4255.5
$ws.Range("H137").Value2 = 3660.3225
$ws.Range("I137").Value2 = 2051.6316
$ws.Range("J137").Value2 = 6207.4165
$ws.Range("K137").Value2 = 6154.8948
$ws.Range("L137").Value2 = 18622.2495
$ws.Range("M137").Value2 = -3604.8948
$ws.Range("N137").Value2 = -23722.2495
$ws.Range("H138").Value2 = 4073.2896
$ws.Range("I138").Value2 = 1841
$ws.Range("J138").Value2 = 5880.381
$ws.Range("K138").Value2 = 5523
$ws.Range("L138").Value2 = 17641.143
$ws.Range("M138").Value2 = -383
$ws.Range("N138").Value2 = -27921.143
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 731.25
$ws.Range("I2").Value2 = 642.5
$ws.Range("J2").Value2 = 820
$ws.Range("K2").Value2 = 642.5
$ws.Range("L2").Value2 = 820
$ws.Range("M2").Value2 = -529.5
$ws.Range("N2").Value2 = -1046
$ws.Range("H32").Value2 = 4067.575
$ws.Range("I32").Value2 = 3542.0557
$ws.Range("K32").Value2 = 3542.0557
$ws.Range("M32").Value2 = -3255.0557
$ws.Range("H61").Value2 = 1019.2
$ws.Range("I61").Value2 = 1019.2
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 1019.2
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -807.2
$ws.Range("H74").Value2 = 2010.8684
$ws.Range("I74").Value2 = 2005.7838
$ws.Range("K74").Value2 = 2005.7838
$ws.Range("M74").Value2 = -1131.7838
$ws.Range("H77").Value2 = 2010.8684
$ws.Range("I77").Value2 = 2005.7838
$ws.Range("K77").Value2 = 10028.919
$ws.Range("M77").Value2 = -5660.919
$ws.Range("H110").Value2 = 1555.5
$ws.Range("I110").Value2 = 2000
$ws.Range("J110").Value2 = 1111
$ws.Range("K110").Value2 = 2000
$ws.Range("L110").Value2 = 1111
$ws.Range("M110").Value2 = 45
$ws.Range("N110").Value2 = -5201
$ws.Range("H116").Value2 = 731.25
$ws.Range("I116").Value2 = 642.5
$ws.Range("J116").Value2 = 820
$ws.Range("K116").Value2 = 642.5
$ws.Range("L116").Value2 = 820
$ws.Range("M116").Value2 = 1651.5
$ws.Range("N116").Value2 = -5408
$ws.Range("H132").Value2 = 2881.9048
$ws.Range("I132").Value2 = 2477.7646
$ws.Range("K132").Value2 = 7433.293799999999
$ws.Range("M132").Value2 = -4903.293799999999
$ws.Range("H136").Value2 = 1019.2
$ws.Range("I136").Value2 = 1019.2
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 3057.6
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -507.6000000000004
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 731.25
$ws.Range("I3").Value2 = 642.5
$ws.Range("J3").Value2 = 820
$ws.Range("K3").Value2 = 642.5
$ws.Range("L3").Value2 = 820
$ws.Range("M3").Value2 = -528.5
$ws.Range("N3").Value2 = -1048
$ws.Range("H105").Value2 = 3015.1
$ws.Range("I105").Value2 = 3016.7778
$ws.Range("J105").Value2 = 3000
$ws.Range("K105").Value2 = 3016.7778
$ws.Range("L105").Value2 = 3000
$ws.Range("M105").Value2 = -1269.7778
$ws.Range("N105").Value2 = -6494
$ws.Range("H107").Value2 = 2123.923
$ws.Range("I107").Value2 = 1812.1
$ws.Range("J107").Value2 = 3163.3333
$ws.Range("K107").Value2 = 1812.1
$ws.Range("L107").Value2 = 3163.3333
$ws.Range("M107").Value2 = 107.9000000000001
$ws.Range("N107").Value2 = -7003.3333
$ws.Range("H132").Value2 = 109866.336
$ws.Range("I132").Value2 = 100000
$ws.Range("J132").Value2 = 114799.5
$ws.Range("K132").Value2 = 100000
$ws.Range("L132").Value2 = 114799.5
$ws.Range("M132").Value2 = -94940
$ws.Range("N132").Value2 = -124919.5
$ws.Range("H134").Value2 = 3490.7334
$ws.Range("I134").Value2 = 3554.3572
$ws.Range("K134").Value2 = 10663.0716
$ws.Range("M134").Value2 = -8128.071599999999
$ws.Range("H135").Value2 = 53836.25
$ws.Range("J135").Value2 = 53836.25
$ws.Range("L135").Value2 = 53836.25
$ws.Range("N135").Value2 = -63976.25
$ws.Range("H138").Value2 = 122199.5
$ws.Range("J138").Value2 = 122199.5
$ws.Range("L138").Value2 = 122199.5
$ws.Range("N138").Value2 = -132479.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 1
$ws.Range("J4").Value2 = 1
$ws.Range("L4").Value2 = 1
$ws.Range("N4").Value2 = -225
$ws.Range("H31").Value2 = 2927
$ws.Range("I31").Value2 = 2967.5715
$ws.Range("K31").Value2 = 2967.5715
$ws.Range("M31").Value2 = -2672.5715
$ws.Range("H34").Value2 = 2927
$ws.Range("I34").Value2 = 2967.5715
$ws.Range("K34").Value2 = 2967.5715
$ws.Range("M34").Value2 = -2765.5715
$ws.Range("H132").Value2 = 2608.818
$ws.Range("I132").Value2 = 2568.3
$ws.Range("J132").Value2 = 3014
$ws.Range("K132").Value2 = 7704.900000000001
$ws.Range("L132").Value2 = 9042
$ws.Range("M132").Value2 = -5174.900000000001
$ws.Range("N132").Value2 = -14102
$ws.Range("H134").Value2 = 2413.8125
$ws.Range("I134").Value2 = 2434.7334
$ws.Range("K134").Value2 = 7304.2002
$ws.Range("M134").Value2 = -4769.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("H92").Value2 = 2000
$ws.Range("I92").Value2 = 2000
$ws.Range("J92").Value2 = 2000
$ws.Range("K92").Value2 = 6000
$ws.Range("L92").Value2 = 6000
$ws.Range("M92").Value2 = -4752
$ws.Range("N92").Value2 = -8496
$ws.Range("H137").Value2 = 6071.2856
$ws.Range("I137").Value2 = 11616.667
$ws.Range("J137").Value2 = 1912.25
$ws.Range("K137").Value2 = 34850.001
$ws.Range("L137").Value2 = 5736.75
$ws.Range("M137").Value2 = -29750.001
$ws.Range("N137").Value2 = -15936.75
$ws.Range("H141").Value2 = 3749.5
$ws.Range("I141").Value2 = 3749.5
$ws.Range("K141").Value2 = 11248.5
$ws.Range("M141").Value2 = -6068.5
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 2403.25
$ws.Range("I93").Value2 = 2403.25
$ws.Range("K93").Value2 = 2403.25
$ws.Range("M93").Value2 = -1155.25
$ws.Range("H132").Value2 = 3075.0588
$ws.Range("I132").Value2 = 2213.8333
$ws.Range("J132").Value2 = 3544.818
$ws.Range("K132").Value2 = 6641.499899999999
$ws.Range("L132").Value2 = 10634.454
$ws.Range("M132").Value2 = -4111.499899999999
$ws.Range("N132").Value2 = -15694.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 2380.0571
$ws.Range("I132").Value2 = 1644.25
$ws.Range("K132").Value2 = 4932.75
$ws.Range("M132").Value2 = -2402.75
$ws.Range("H136").Value2 = 999
$ws.Range("I136").Value2 = 919.16
$ws.Range("K136").Value2 = 2757.48
$ws.Range("M136").Value2 = -207.48
